# Fruta / hortaliza, semanal
# Inserts a new weekly record at row 77 (Macroferia Regional de Talca - Haba),
# shifting the existing rows 77:107 down to 78:108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data (rows 77..107) down one row to make room for the new record.
$ws.Rows("77").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(77, 1).Value = 5
$ws.Cells.Item(77, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(77, 3).Value = "Maule"
$ws.Cells.Item(77, 4).Value = 44876
$ws.Cells.Item(77, 5).Value = 7
$ws.Cells.Item(77, 6).Value = 100112026
$ws.Cells.Item(77, 7).Value = "Haba"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 400
$ws.Cells.Item(77, 11).Value = 8000
$ws.Cells.Item(77, 12).Value = 8000
$ws.Cells.Item(77, 13).Value = 8000
$ws.Cells.Item(77, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(77, 15).Value = "Región del Maule"
$ws.Cells.Item(77, 16).Value = 320
$ws.Cells.Item(77, 17).Value = 25
$ws.Cells.Item(77, 18).Value = "Hortaliza"
